$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 64, shifting existing rows 64:133 down to 65:134.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with this week's data (mirrors the
# columns used throughout the rest of the table for this seller/product).
$ws.Range("A64").Value = 1
$ws.Range("B64").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C64").Value = "Arica y Parinacota"
$ws.Range("D64").Value = 45225
$ws.Range("E64").Value = 15
$ws.Range("F64").Value = 100112040
$ws.Range("G64").Value = "Cilantro"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 250
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 2000
$ws.Range("M64").Value = 1750
$ws.Range("N64").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O64").Value = "Región de Arica y Parinacota"
$ws.Range("P64").Value = 1167
$ws.Range("Q64").Value = 1.5
$ws.Range("R64").Value = "Hortaliza"
